$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price cells being updated so numeric-looking
# strings (e.g. "227.40") are stored verbatim as text, matching the
# original inline-string cell contents (no precision/zero loss). Only the
# touched D-column cells get this so untouched cells keep their style.
# (Comma-separated multi-area Range refs aren't honored, so loop instead.)
$dRows = @(2,3,5,6,7,9,12,13,14,15,16,17,18,19,20,21,22,26,27,29,32,33,34,36,37,38,40,41,42,43,45,47,48,49,50,51)
foreach ($r in $dRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

$ws.Range("D2").Value = '37.719.20'
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").Value = '2.035.15'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '227.40'
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("D6").Value = '0.607'
$ws.Range("E6").Value = '  -0.58%  '
$ws.Range("D7").Value = '60.04'
$ws.Range("E7").Value = '  +0.50%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '0.376'
$ws.Range("E9").Value = '  -1.02%  '
$ws.Range("E10").Value = '  +1.93%  '
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("D12").Value = '2.332.27'
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("D13").Value = '14.53'
$ws.Range("E13").Value = '  +0.08%  '
$ws.Range("D14").Value = '21.10'
$ws.Range("E14").Value = '  +0.99%  '
$ws.Range("D15").Value = '0.774'
$ws.Range("E15").Value = '  +3.25%  '
$ws.Range("D16").Value = '5.32'
$ws.Range("E16").Value = '  +1.96%  '
$ws.Range("D17").Value = '2.025.25'
$ws.Range("E17").Value = '  -0.80%  '
$ws.Range("D18").Value = '37.653.84'
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("D19").Value = '5.96'
$ws.Range("E19").Value = '  -2.20%  '
$ws.Range("D20").Value = '69.40'
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("D21").Value = '0.0₃0821'
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("D22").Value = '223.94'
$ws.Range("E22").Value = '  +0.27%  '
$ws.Range("E23").Value = '  +0.14%  '
$ws.Range("E24").Value = '  -0.17%  '
$ws.Range("E25").Value = '  +4.22%  '
$ws.Range("D26").Value = '167.80'
$ws.Range("E26").Value = '  +1.27%  '
$ws.Range("D27").Value = '9.35'
$ws.Range("E27").Value = '  +2.76%  '
$ws.Range("E28").Value = '  -0.52%  '
$ws.Range("D29").Value = '18.77'
$ws.Range("E29").Value = '  -0.26%  '
$ws.Range("E30").Value = '  -1.61%  '
$ws.Range("E31").Value = '  +0.99%  '
$ws.Range("D32").Value = '2.21'
$ws.Range("E32").Value = '  +8.36%  '
$ws.Range("D33").Value = '4.37'
$ws.Range("E33").Value = '  -1.19%  '
$ws.Range("D34").Value = '0.0605'
$ws.Range("E34").Value = '  +0.53%  '
$ws.Range("E35").Value = '  -0.12%  '
$ws.Range("D36").Value = '6.52'
$ws.Range("E36").Value = '  +4.73%  '
$ws.Range("D37").Value = '2.35'
$ws.Range("E37").Value = '  +4.37%  '
$ws.Range("D38").Value = '3.39'
$ws.Range("E38").Value = '  +6.14%  '
$ws.Range("E39").Value = '  +0.06%  '
$ws.Range("D40").Value = '17.90'
$ws.Range("E40").Value = '  +9.08%  '
$ws.Range("D41").Value = '1.527.34'
$ws.Range("E41").Value = '  -0.14%  '
$ws.Range("D42").Value = '97.03'
$ws.Range("E42").Value = '  +0.80%  '
$ws.Range("D43").Value = '0.0215'
$ws.Range("E43").Value = '  -0.34%  '
$ws.Range("E44").Value = '  +0.55%  '
$ws.Range("D45").Value = '0.0908'
$ws.Range("E45").Value = '  -0.79%  '
$ws.Range("E46").Value = '  +3.64%  '
$ws.Range("D47").Value = '1.11'
$ws.Range("E47").Value = '  +0.07%  '
$ws.Range("D48").Value = '1.01'
$ws.Range("E48").Value = '  +0.87%  '
$ws.Range("B49").Value = 'MXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D49").Value = '2.94'
$ws.Range("E49").Value = '  -0.64%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").Value = '7.05'
$ws.Range("E50").Value = '  +0.48%  '
$ws.Range("D51").Value = '2.224.40'
$ws.Range("E51").Value = '  +0.24%  '
